$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 412.7143
$ws.Range("I12").Value = 355.8
$ws.Range("K12").Value = 355.8
$ws.Range("M12").Value = -185.8
$ws.Range("H19").Value = 2257.6
$ws.Range("J19").Value = 2270.8
$ws.Range("L19").Value = 2270.8
$ws.Range("N19").Value = -2620.8
$ws.Range("H33").Value = 124.8
$ws.Range("I33").Value = 125.4
$ws.Range("K33").Value = 125.4
$ws.Range("M33").Value = 103.6
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H42").Value = 2555.5
$ws.Range("J42").Value = 3380.6667
$ws.Range("L42").Value = 10142.0001
$ws.Range("N42").Value = -10602.0001
$ws.Range("H70").Value = 3499.5833
$ws.Range("J70").Value = 3499.5833
$ws.Range("L70").Value = 10498.7499
$ws.Range("N70").Value = -11038.7499
$ws.Range("H73").Value = 3499.5833
$ws.Range("J73").Value = 3499.5833
$ws.Range("L73").Value = 10498.7499
$ws.Range("N73").Value = -12370.7499
$ws.Range("H98").Value = 1764.1428
$ws.Range("I98").Value = 1568.6
$ws.Range("K98").Value = 1568.6
$ws.Range("M98").Value = -70.59999999999991
$ws.Range("H107").Value = 176
$ws.Range("J107").Value = 166
$ws.Range("L107").Value = 166
$ws.Range("N107").Value = -4006
$ws.Range("H118").Value = 377.83334
$ws.Range("I118").Value = 377.83334
$ws.Range("K118").Value = 1133.50002
$ws.Range("M118").Value = 523.4999800000001
$ws.Range("H122").Value = 1764.1428
$ws.Range("I122").Value = 1568.6
$ws.Range("K122").Value = 4705.799999999999
$ws.Range("M122").Value = -2255.799999999999
$ws.Range("H138").Value = 3504.5151
$ws.Range("I138").Value = 2499.25
$ws.Range("K138").Value = 7497.75
$ws.Range("M138").Value = -2357.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2008.8889
$ws.Range("I45").Value = 1776
$ws.Range("J45").Value = 2300
$ws.Range("K45").Value = 1776
$ws.Range("L45").Value = 2300
$ws.Range("M45").Value = -1399
$ws.Range("N45").Value = -3054
$ws.Range("H74").Value = 3689.5833
$ws.Range("I74").Value = 3689.5833
$ws.Range("K74").Value = 3689.5833
$ws.Range("M74").Value = -2815.5833
$ws.Range("H77").Value = 3689.5833
$ws.Range("I77").Value = 3689.5833
$ws.Range("K77").Value = 18447.9165
$ws.Range("M77").Value = -14079.9165
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 1244.7241
$ws.Range("I132").Value = 1123.92
$ws.Range("K132").Value = 3371.76
$ws.Range("M132").Value = -841.7600000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 696.75
$ws.Range("I80").Value = 688.3333
$ws.Range("J80").Value = 701.8
$ws.Range("K80").Value = 688.3333
$ws.Range("L80").Value = 701.8
$ws.Range("M80").Value = 309.6667
$ws.Range("N80").Value = -2697.8
$ws.Range("H83").Value = 696.75
$ws.Range("I83").Value = 688.3333
$ws.Range("J83").Value = 701.8
$ws.Range("K83").Value = 3472.5
$ws.Range("L83").Value = 3517
$ws.Range("M83").Value = 1550.3335
$ws.Range("N83").Value = -13493
$ws.Range("H96").Value = 20000
$ws.Range("I96").Value = 20000
$ws.Range("K96").Value = 20000
$ws.Range("M96").Value = -17254
$ws.Range("H105").Value = 10000
$ws.Range("I105").Value = 9000
$ws.Range("J105").Value = 11000
$ws.Range("K105").Value = 9000
$ws.Range("L105").Value = 11000
$ws.Range("M105").Value = -7253
$ws.Range("N105").Value = -14494
$ws.Range("H134").Value = 2476.9333
$ws.Range("I134").Value = 2476.9333
$ws.Range("K134").Value = 7430.7999
$ws.Range("M134").Value = -4895.7999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 177.8
$ws.Range("I7").Value = 113
$ws.Range("K7").Value = 113
$ws.Range("M7").Value = 0
$ws.Range("H16").Value = 1423.2
$ws.Range("I16").Value = 1489.8572
$ws.Range("J16").Value = 1267.6666
$ws.Range("K16").Value = 1489.8572
$ws.Range("L16").Value = 1267.6666
$ws.Range("M16").Value = -1202.8572
$ws.Range("N16").Value = -1841.6666
$ws.Range("H58").Value = 5548.75
$ws.Range("I58").Value = 5548.75
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 5548.75
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -5345.75
$ws.Range("N58").ClearContents()
$ws.Range("H103").Value = 3262
$ws.Range("I103").Value = 3262
$ws.Range("K103").Value = 3262
$ws.Range("M103").Value = -2090
$ws.Range("H113").Value = 1423.2
$ws.Range("I113").Value = 1489.8572
$ws.Range("J113").Value = 1267.6666
$ws.Range("K113").Value = 1489.8572
$ws.Range("L113").Value = 1267.6666
$ws.Range("M113").Value = 680.1428000000001
$ws.Range("N113").Value = -5607.6666
$ws.Range("H122").Value = 917.8
$ws.Range("I122").Value = 917.8
$ws.Range("K122").Value = 2753.4
$ws.Range("M122").Value = -303.3999999999996
$ws.Range("H132").Value = 3098.3333
$ws.Range("I132").Value = 3038
$ws.Range("J132").Value = 3400
$ws.Range("K132").Value = 9114
$ws.Range("L132").Value = 10200
$ws.Range("M132").Value = -6584
$ws.Range("N132").Value = -15260
$ws.Range("H136").Value = 5548.75
$ws.Range("I136").Value = 5548.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 16646.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -14096.25
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 75.77778000000001
$ws.Range("J38").Value = 136.33333
$ws.Range("L38").Value = 408.99999
$ws.Range("N38").Value = -1102.99999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3662
$ws.Range("I102").Value = 3662
$ws.Range("K102").Value = 3662
$ws.Range("M102").Value = -2040
$ws.Range("H113").Value = 2580.75
$ws.Range("I113").Value = 2729.4
$ws.Range("J113").Value = 2333
$ws.Range("K113").Value = 2729.4
$ws.Range("L113").Value = 2333
$ws.Range("M113").Value = -559.4000000000001
$ws.Range("N113").Value = -6673
$ws.Range("H132").Value = 1700
$ws.Range("I132").Value = 1700
$ws.Range("K132").Value = 5100
$ws.Range("M132").Value = -2570

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 39998
$ws.Range("I74").Value = 20000
$ws.Range("K74").Value = 20000
$ws.Range("M74").Value = -19002
$ws.Range("H77").Value = 39998
$ws.Range("I77").Value = 20000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55008
$ws.Range("H100").Value = 3304.5
$ws.Range("I100").Value = 3304.5
$ws.Range("K100").Value = 3304.5
$ws.Range("M100").Value = -2763.5
$ws.Range("H132").Value = 7282
$ws.Range("I132").Value = 6600.6665
$ws.Range("J132").Value = 8048.5
$ws.Range("K132").Value = 19801.9995
$ws.Range("L132").Value = 24145.5
$ws.Range("M132").Value = -17271.9995
$ws.Range("N132").Value = -29205.5
$ws.Range("H136").Value = 3494
$ws.Range("I136").Value = 3494
$ws.Range("K136").Value = 10482
$ws.Range("M136").Value = -7932

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 870
$ws.Range("J107").Value = 868.75
$ws.Range("L107").Value = 2606.25
$ws.Range("N107").Value = -6446.25
$ws.Range("H122").Value = 2944
$ws.Range("I122").Value = 2842.111
$ws.Range("K122").Value = 8526.332999999999
$ws.Range("M122").Value = -6076.332999999999
$ws.Range("H126").Value = 1499.6666
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1499.6666
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 4498.9998
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -9438.9998
$ws.Range("H132").Value = 722.35297
$ws.Range("I132").Value = 618.73334
$ws.Range("K132").Value = 1856.20002
$ws.Range("M132").Value = 673.79998
